$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.556.06"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "1.565.73"
$ws.Range("E3").Value = "  -1.70%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "211.74"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -1.43%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.495"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("E7").Value = "  +0.16%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "46.30"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +5.48%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "24.18"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("E12").Value = "  -0.74%  "
$ws.Range("D13").Value = "1.788.79"
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("D14").Value = "1.560.30"
$ws.Range("E14").Value = "  -2.19%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.520"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -2.15%  "
$ws.Range("D16").Value = "28.538.73"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("E17").Value = "  -3.12%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "62.00"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -3.20%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "227.31"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -2.60%  "
$ws.Range("E20").Value = "  -2.37%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "7.32"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("E23").Value = "  -6.92%  "
$ws.Range("E24").Value = "  -3.36%  "
$ws.Range("E25").Value = "  +5.81%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "150.72"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("E27").Value = "  -2.56%  "
$ws.Range("E28").Value = "  -3.07%  "
$ws.Range("E29").Value = "  -3.75%  "
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("E31").Value = "  -1.99%  "
$ws.Range("E32").Value = "  -3.87%  "
$ws.Range("E33").Value = "  -1.89%  "
$ws.Range("E34").Value = "  -0.26%  "
$ws.Range("D35").Value = "1.395.65"
$ws.Range("E35").Value = "  -1.67%  "
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("E37").Value = "  -3.99%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "2.36"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +1.19%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.58"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +0.71%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.536"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -1.56%  "
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("E43").Value = "  -3.59%  "
$ws.Range("E44").Value = "  +1.47%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "5.51"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -3.91%  "
$ws.Range("E46").Value = "  -0.13%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "62.63"
$cell.Style = "Normal"
$ws.Range("D48").Value = "1.701.63"
$ws.Range("E48").Value = "  -1.70%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "86.09"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -2.17%  "
$ws.Range("D50").Value = "0.0₆0103"
$ws.Range("E50").Value = "  -0.40%  "
$ws.Range("E51").Value = "  -0.96%  "
